$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "sure bitti"
$ws.Range("A5").Select()
